# Auto-generated script applying the price/profit data refresh
# captured by the scheduled Sheets runner (see commit message).
$wb = $excel.ActiveWorkbook

# ----- ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 1069.8   # H6: was 1250.25
$ws.Cells.Item(6, 9).Value = 586.75   # I6: was 666.3333
$ws.Cells.Item(6, 11).Value = 1760.25   # K6: was 1998.9999
$ws.Cells.Item(6, 13).Value = -1648.25   # M6: was -1886.9999
$ws.Cells.Item(40, 8).Value = 3104.9375   # H40: was 3095.5625
$ws.Cells.Item(40, 9).Value = 2306.5833   # I40: was 2194.5386
$ws.Cells.Item(40, 10).Value = 5500   # J40: was 7000
$ws.Cells.Item(40, 11).Value = 2306.5833   # K40: was 2194.5386
$ws.Cells.Item(40, 12).Value = 5500   # L40: was 7000
$ws.Cells.Item(40, 13).Value = -2131.5833   # M40: was -2019.5386
$ws.Cells.Item(40, 14).Value = -5850   # N40: was -7350
$ws.Cells.Item(62, 8).Value = 8212   # H62: was 8212.143
$ws.Cells.Item(62, 9).Value = 8214.166999999999   # I62: was 8214.333000000001
$ws.Cells.Item(62, 11).Value = 8214.166999999999   # K62: was 8214.333000000001
$ws.Cells.Item(62, 13).Value = -7590.166999999999   # M62: was -7590.333000000001
$ws.Cells.Item(65, 8).Value = 8212   # H65: was 8212.143
$ws.Cells.Item(65, 9).Value = 8214.166999999999   # I65: was 8214.333000000001
$ws.Cells.Item(65, 11).Value = 41070.835   # K65: was 41071.665
$ws.Cells.Item(65, 13).Value = -37950.835   # M65: was -37951.665
$ws.Cells.Item(70, 8).Value = 1880.2727   # H70: was 1880.7273
$ws.Cells.Item(70, 9).Value = 1846.5   # I70: was 1849
$ws.Cells.Item(70, 11).Value = 5539.5   # K70: was 5547
$ws.Cells.Item(70, 13).Value = -5269.5   # M70: was -5277
$ws.Cells.Item(73, 8).Value = 1880.2727   # H73: was 1880.7273
$ws.Cells.Item(73, 9).Value = 1846.5   # I73: was 1849
$ws.Cells.Item(73, 11).Value = 5539.5   # K73: was 5547
$ws.Cells.Item(73, 13).Value = -4603.5   # M73: was -4611
$ws.Cells.Item(86, 8).Value = 9000   # H86: was 0
$ws.Cells.Item(86, 10).Value = 9000   # J86: was 0
$ws.Cells.Item(86, 12).Value = 9000   # L86: was 0
$ws.Cells.Item(86, 14).Value = -11246   # N86: was None
$ws.Cells.Item(89, 8).Value = 9000   # H89: was 0
$ws.Cells.Item(89, 10).Value = 9000   # J89: was 0
$ws.Cells.Item(89, 12).Value = 45000   # L89: was 0
$ws.Cells.Item(89, 14).Value = -56232   # N89: was None
$ws.Cells.Item(100, 8).Value = 4326.5557   # H100: was 3744
$ws.Cells.Item(100, 9).Value = 4742.375   # I100: was 3993.4546
$ws.Cells.Item(100, 11).Value = 4742.375   # K100: was 3993.4546
$ws.Cells.Item(100, 13).Value = -4201.375   # M100: was -3452.4546
$ws.Cells.Item(116, 8).Value = 4494.6   # H116: was 4496.6
$ws.Cells.Item(116, 9).Value = 4497   # I116: was 4496
$ws.Cells.Item(116, 10).Value = 4491   # J116: was 4499
$ws.Cells.Item(116, 11).Value = 4497   # K116: was 4496
$ws.Cells.Item(116, 12).Value = 4491   # L116: was 4499
$ws.Cells.Item(116, 13).Value = -1055   # M116: was -1054
$ws.Cells.Item(116, 14).Value = -11375   # N116: was -11383
$ws.Cells.Item(131, 8).Value = 4140   # H131: was 3566.6667
$ws.Cells.Item(138, 8).Value = 3354.1853   # H138: was 3390.2964
$ws.Cells.Item(138, 10).Value = 3524.5652   # J138: was 3566.9565
$ws.Cells.Item(138, 12).Value = 10573.6956   # L138: was 10700.8695
$ws.Cells.Item(138, 14).Value = -20853.6956   # N138: was -20980.8695

# ----- ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(102, 8).Value = 2612.9092   # H102: was 2307.2144
$ws.Cells.Item(102, 9).Value = 2805.25   # I102: was 2363.7273
$ws.Cells.Item(102, 11).Value = 2805.25   # K102: was 2363.7273
$ws.Cells.Item(102, 13).Value = -1183.25   # M102: was -741.7273

# ----- BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 5622   # H20: was 4836.1816
$ws.Cells.Item(20, 9).Value = 1419.8   # I20: was 1283.1666
$ws.Cells.Item(20, 10).Value = 10874.75   # J20: was 9099.799999999999
$ws.Cells.Item(20, 11).Value = 1419.8   # K20: was 1283.1666
$ws.Cells.Item(20, 12).Value = 10874.75   # L20: was 9099.799999999999
$ws.Cells.Item(20, 13).Value = -1172.8   # M20: was -1036.1666
$ws.Cells.Item(20, 14).Value = -11368.75   # N20: was -9593.799999999999
$ws.Cells.Item(86, 8).Value = 5236   # H86: was 4814
$ws.Cells.Item(86, 9).Value = 4333   # I86: was 2891.5
$ws.Cells.Item(86, 10).Value = 5494   # J86: was 6095.6665
$ws.Cells.Item(86, 11).Value = 4333   # K86: was 2891.5
$ws.Cells.Item(86, 12).Value = 5494   # L86: was 6095.6665
$ws.Cells.Item(86, 13).Value = -3210   # M86: was -1768.5
$ws.Cells.Item(86, 14).Value = -7740   # N86: was -8341.666499999999
$ws.Cells.Item(89, 8).Value = 5236   # H89: was 4814
$ws.Cells.Item(89, 9).Value = 4333   # I89: was 2891.5
$ws.Cells.Item(89, 10).Value = 5494   # J89: was 6095.6665
$ws.Cells.Item(89, 11).Value = 21665   # K89: was 14457.5
$ws.Cells.Item(89, 12).Value = 27470   # L89: was 30478.3325
$ws.Cells.Item(89, 13).Value = -16049   # M89: was -8841.5
$ws.Cells.Item(89, 14).Value = -38702   # N89: was -41710.3325
$ws.Cells.Item(99, 8).Value = 1076   # H99: was 882.1667
$ws.Cells.Item(99, 9).Value = 1076   # I99: was 882.1667
$ws.Cells.Item(99, 11).Value = 1076   # K99: was 882.1667
$ws.Cells.Item(99, 13).Value = 422   # M99: was 615.8333
$ws.Cells.Item(100, 8).Value = 42773.668   # H100: was 43806
$ws.Cells.Item(100, 10).Value = 42773.668   # J100: was 43806
$ws.Cells.Item(100, 12).Value = 42773.668   # L100: was 43806
$ws.Cells.Item(100, 14).Value = -44937.668   # N100: was -45970

# ----- CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(29, 8).Value = 0   # H29: was 50000
$ws.Cells.Item(29, 10).Value = 0   # J29: was 50000
$ws.Cells.Item(29, 12).ClearContents()   # L29: was 50000
$ws.Cells.Item(29, 14).Value = 0   # N29: was -50586
$ws.Cells.Item(31, 8).Value = 2928.875   # H31: was 2990.75
$ws.Cells.Item(31, 9).Value = 2772   # I31: was 2787.8333
$ws.Cells.Item(31, 10).Value = 3399.5   # J31: was 3599.5
$ws.Cells.Item(31, 11).Value = 2772   # K31: was 2787.8333
$ws.Cells.Item(31, 12).Value = 3399.5   # L31: was 3599.5
$ws.Cells.Item(31, 13).Value = -2477   # M31: was -2492.8333
$ws.Cells.Item(31, 14).Value = -3989.5   # N31: was -4189.5
$ws.Cells.Item(34, 8).Value = 2928.875   # H34: was 2990.75
$ws.Cells.Item(34, 9).Value = 2772   # I34: was 2787.8333
$ws.Cells.Item(34, 10).Value = 3399.5   # J34: was 3599.5
$ws.Cells.Item(34, 11).Value = 2772   # K34: was 2787.8333
$ws.Cells.Item(34, 12).Value = 3399.5   # L34: was 3599.5
$ws.Cells.Item(34, 13).Value = -2570   # M34: was -2585.8333
$ws.Cells.Item(34, 14).Value = -3803.5   # N34: was -4003.5
$ws.Cells.Item(50, 8).Value = 47028   # H50: was 45992
$ws.Cells.Item(50, 9).Value = 50000   # I50: was 0
$ws.Cells.Item(50, 10).Value = 45542   # J50: was 45992
$ws.Cells.Item(50, 11).Value = 50000   # K50: was 0
$ws.Cells.Item(50, 12).Value = 45542   # L50: was 45992
$ws.Cells.Item(50, 13).Value = -49375   # M50: was None
$ws.Cells.Item(50, 14).Value = -46792   # N50: was -47242
$ws.Cells.Item(86, 8).Value = 7070.4287   # H86: was 7070.5713
$ws.Cells.Item(86, 9).Value = 7374.75   # I86: was 6499.8
$ws.Cells.Item(86, 10).Value = 6664.6665   # J86: was 8497.5
$ws.Cells.Item(86, 11).Value = 7374.75   # K86: was 6499.8
$ws.Cells.Item(86, 12).Value = 6664.6665   # L86: was 8497.5
$ws.Cells.Item(86, 13).Value = -6251.75   # M86: was -5376.8
$ws.Cells.Item(86, 14).Value = -8910.666499999999   # N86: was -10743.5
$ws.Cells.Item(89, 8).Value = 7070.4287   # H89: was 7070.5713
$ws.Cells.Item(89, 9).Value = 7374.75   # I89: was 6499.8
$ws.Cells.Item(89, 10).Value = 6664.6665   # J89: was 8497.5
$ws.Cells.Item(89, 11).Value = 36873.75   # K89: was 32499
$ws.Cells.Item(89, 12).Value = 33323.3325   # L89: was 42487.5
$ws.Cells.Item(89, 13).Value = -31257.75   # M89: was -26883
$ws.Cells.Item(89, 14).Value = -44555.3325   # N89: was -53719.5
$ws.Cells.Item(107, 8).Value = 1874.6666   # H107: was 1791.6154
$ws.Cells.Item(107, 10).Value = 866.3333   # J107: was 848.5
$ws.Cells.Item(107, 12).Value = 866.3333   # L107: was 848.5
$ws.Cells.Item(107, 14).Value = -4706.3333   # N107: was -4688.5

# ----- CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(6, 8).Value = 184.6   # H6: was 176.125
$ws.Cells.Item(6, 9).Value = 184.6   # I6: was 176.125
$ws.Cells.Item(6, 11).Value = 553.8   # K6: was 528.375
$ws.Cells.Item(6, 13).Value = -440.8   # M6: was -415.375
$ws.Cells.Item(21, 8).Value = 1000   # H21: was 0
$ws.Cells.Item(21, 9).Value = 1000   # I21: was 0
$ws.Cells.Item(21, 11).Value = 3000   # K21: was 0
$ws.Cells.Item(21, 13).Value = -2827   # M21: was None
$ws.Cells.Item(42, 8).Value = 5980   # H42: was 254485
$ws.Cells.Item(42, 10).Value = 5980   # J42: was 254485
$ws.Cells.Item(42, 12).Value = 17940   # L42: was 763455
$ws.Cells.Item(42, 14).Value = -19008   # N42: was -764523
$ws.Cells.Item(50, 8).Value = 418.33334   # H50: was 368.33334
$ws.Cells.Item(50, 9).Value = 502.5   # I50: was 368.33334
$ws.Cells.Item(50, 10).Value = 250   # J50: was 0
$ws.Cells.Item(50, 11).Value = 1507.5   # K50: was 1105.00002
$ws.Cells.Item(50, 12).Value = 750   # L50: was 0
$ws.Cells.Item(50, 13).Value = -1026.5   # M50: was -624.0000199999999
$ws.Cells.Item(50, 14).Value = -1712   # N50: was None
$ws.Cells.Item(53, 8).Value = 418.33334   # H53: was 368.33334
$ws.Cells.Item(53, 9).Value = 502.5   # I53: was 368.33334
$ws.Cells.Item(53, 10).Value = 250   # J53: was 0
$ws.Cells.Item(53, 11).Value = 1507.5   # K53: was 1105.00002
$ws.Cells.Item(53, 12).Value = 750   # L53: was 0
$ws.Cells.Item(53, 13).Value = -1026.5   # M53: was -624.0000199999999
$ws.Cells.Item(53, 14).Value = -1712   # N53: was None
$ws.Cells.Item(113, 8).Value = 1544.7273   # H113: was 1599.5
$ws.Cells.Item(113, 9).Value = 873.5   # I113: was 832.3333
$ws.Cells.Item(113, 11).Value = 2620.5   # K113: was 2496.9999
$ws.Cells.Item(113, 13).Value = -450.5   # M113: was -326.9998999999998

# ----- GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(32, 8).Value = 0   # H32: was 40000
$ws.Cells.Item(32, 10).Value = 0   # J32: was 40000
$ws.Cells.Item(32, 12).ClearContents()   # L32: was 40000
$ws.Cells.Item(32, 14).Value = 0   # N32: was -40592
$ws.Cells.Item(42, 8).Value = 0   # H42: was 120000
$ws.Cells.Item(42, 10).Value = 0   # J42: was 120000
$ws.Cells.Item(42, 12).ClearContents()   # L42: was 120000
$ws.Cells.Item(42, 14).Value = 0   # N42: was -120970
$ws.Cells.Item(80, 8).Value = 1633.3334   # H80: was 1621.5555
$ws.Cells.Item(80, 9).Value = 1633.3334   # I80: was 1226.8572
$ws.Cells.Item(80, 10).Value = 0   # J80: was 3003
$ws.Cells.Item(80, 11).Value = 1633.3334   # K80: was 1226.8572
$ws.Cells.Item(80, 12).Value = 0   # L80: was 3003
$ws.Cells.Item(80, 13).ClearContents()   # M80: was -228.8571999999999
$ws.Cells.Item(80, 14).Value = -635.3334   # N80: was -4999
$ws.Cells.Item(83, 8).Value = 1633.3334   # H83: was 1621.5555
$ws.Cells.Item(83, 9).Value = 1633.3334   # I83: was 1226.8572
$ws.Cells.Item(83, 10).Value = 0   # J83: was 3003
$ws.Cells.Item(83, 11).Value = 8166.666999999999   # K83: was 6134.286
$ws.Cells.Item(83, 12).Value = 0   # L83: was 15015
$ws.Cells.Item(83, 13).ClearContents()   # M83: was -1142.286
$ws.Cells.Item(83, 14).Value = -3174.666999999999   # N83: was -24999
$ws.Cells.Item(103, 8).Value = 0   # H103: was 40000
$ws.Cells.Item(103, 10).Value = 0   # J103: was 40000
$ws.Cells.Item(103, 12).ClearContents()   # L103: was 40000
$ws.Cells.Item(103, 14).Value = 0   # N103: was -42344
$ws.Cells.Item(115, 8).Value = 0   # H115: was 120000
$ws.Cells.Item(115, 10).Value = 0   # J115: was 120000
$ws.Cells.Item(115, 12).ClearContents()   # L115: was 120000
$ws.Cells.Item(115, 14).Value = 0   # N115: was -122350
$ws.Cells.Item(141, 8).Value = 37500   # H141: was 69000
$ws.Cells.Item(141, 10).Value = 37500   # J141: was 69000
$ws.Cells.Item(141, 12).Value = 37500   # L141: was 69000
$ws.Cells.Item(141, 14).Value = -47860   # N141: was -79360

# ----- LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 6030.737   # H22: was 6033.3687
$ws.Cells.Item(22, 10).Value = 8837.5   # J22: was 8843.75
$ws.Cells.Item(22, 12).Value = 8837.5   # L22: was 8843.75
$ws.Cells.Item(22, 14).Value = -9427.5   # N22: was -9433.75
$ws.Cells.Item(27, 8).Value = 6030.737   # H27: was 6033.3687
$ws.Cells.Item(27, 10).Value = 8837.5   # J27: was 8843.75
$ws.Cells.Item(27, 12).Value = 8837.5   # L27: was 8843.75
$ws.Cells.Item(27, 14).Value = -9051.5   # N27: was -9057.75
$ws.Cells.Item(40, 8).Value = 4113   # H40: was 4163.5625
$ws.Cells.Item(40, 9).Value = 4044   # I40: was 4111.4165
$ws.Cells.Item(40, 11).Value = 4044   # K40: was 4111.4165
$ws.Cells.Item(40, 13).Value = -3908   # M40: was -3975.4165
$ws.Cells.Item(46, 8).Value = 7024.5625   # H46: was 8539.923000000001
$ws.Cells.Item(46, 10).Value = 15805.5   # J46: was 31153
$ws.Cells.Item(46, 12).Value = 15805.5   # L46: was 31153
$ws.Cells.Item(46, 14).Value = -16181.5   # N46: was -31529

# ----- WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 1062   # H107: was 996.7778
$ws.Cells.Item(107, 9).Value = 899.4   # I107: was 828.6667
$ws.Cells.Item(107, 11).Value = 2698.2   # K107: was 2486.0001
$ws.Cells.Item(107, 13).Value = -778.1999999999998   # M107: was -566.0001000000002
